# fix(FN-3460): fix invalid facility utilisation values -- all need to
# match as same facility id for all rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: facility utilisation values
$ws.Range("E5").Value = 600000
$ws.Range("F5").Value = 761579.37
$ws.Range("G5").Value = 3938753.8

# Row 6: facility utilisation values
$ws.Range("E6").Value = 600000
$ws.Range("G6").Value = 761579.37

# Column G now lines up with the width already used by columns E:F
# (16.33203125 / bestFit) - only G needs to move, E:F are untouched
$ws.Range("G1").EntireColumn.ColumnWidth = 15.417

# Reflect the edited range in the sheet's active selection
$ws.Range("E5:H6").Select()
